$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing it to stay text, even when the
# string looks numeric (e.g. "1.003" or "27.936.40" would otherwise be
# parsed into a number/date by Excel), then restore the cell formatting
# so no stray number format is left on the cell afterwards.
function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue 'D2' '27.936.40'
$ws.Range("E2").Value = '  +1.52%  '
Set-TextValue 'D3' '1.754.19'
$ws.Range("E3").Value = '  -0.42%  '
Set-TextValue 'D4' '1.003'
$ws.Range("E4").Value = '  -0.59%  '
Set-TextValue 'D5' '336.18'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("E6").Value = '  -0.39%  '
Set-TextValue 'D7' '0.3838'
$ws.Range("E7").Value = '  +0.20%  '
Set-TextValue 'D8' '0.3409'
$ws.Range("E8").Value = '  +0.33%  '
Set-TextValue 'D9' '46.18'
$ws.Range("E9").Value = '  -1.86%  '
Set-TextValue 'D10' '1.120'
$ws.Range("E10").Value = '  -1.42%  '
Set-TextValue 'D11' '0.07239'
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D12' '22.57'
$ws.Range("E12").Value = '  +4.03%  '
$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D13' '1.001'
$ws.Range("E13").Value = '  -0.51%  '
Set-TextValue 'D14' '6.179'
$ws.Range("E14").Value = '  -2.27%  '
Set-TextValue 'D15' '7.150'
$ws.Range("E15").Value = '  +1.84%  '
Set-TextValue 'D16' '1.752.73'
$ws.Range("E16").Value = '  -0.71%  '
Set-TextValue 'D17' '0.00001062'
$ws.Range("E17").Value = '  -0.65%  '
Set-TextValue 'D18' '0.06615'
$ws.Range("E18").Value = '  -0.62%  '
Set-TextValue 'D19' '79.16'
$ws.Range("E19").Value = '  -3.56%  '
Set-TextValue 'D20' '1.000'
$ws.Range("E20").Value = '  -0.52%  '
Set-TextValue 'D21' '16.74'
$ws.Range("E21").Value = '  -3.03%  '
Set-TextValue 'D22' '6.236'
$ws.Range("E22").Value = '  -2.62%  '
Set-TextValue 'D23' '27.931.67'
Set-TextValue 'D24' '11.70'
$ws.Range("E24").Value = '  -2.83%  '
Set-TextValue 'D25' '2.382'
$ws.Range("E25").Value = '  -0.14%  '
Set-TextValue 'D26' '153.92'
$ws.Range("E26").Value = '  +0.68%  '
Set-TextValue 'D27' '19.88'
$ws.Range("E27").Value = '  -3.62%  '
Set-TextValue 'D28' '2.319'
$ws.Range("E28").Value = '  -4.44%  '
Set-TextValue 'D29' '1.952.72'
$ws.Range("E29").Value = '  -0.68%  '
Set-TextValue 'D30' '1.272'
$ws.Range("E30").Value = '  -10.02%  '
Set-TextValue 'D31' '131.85'
$ws.Range("E31").Value = '  -1.97%  '
Set-TextValue 'D32' '4.017'
Set-TextValue 'D33' '5.872'
$ws.Range("E33").Value = '  -3.18%  '
Set-TextValue 'D34' '0.08838'
$ws.Range("E34").Value = '  +1.37%  '
Set-TextValue 'D35' '12.27'
$ws.Range("E35").Value = '  -2.67%  '
Set-TextValue 'D36' '0.6604'
$ws.Range("E36").Value = '  -1.68%  '
Set-TextValue 'D37' '0.02294'
$ws.Range("E37").Value = '  -4.76%  '
Set-TextValue 'D38' '5.163'
$ws.Range("E38").Value = '  -3.15%  '
Set-TextValue 'D39' '0.06183'
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D40' '0.2106'
$ws.Range("E40").Value = '  -3.40%  '
$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D41' '1.500'
$ws.Range("E41").Value = '  -1.72%  '
Set-TextValue 'D42' '1.213'
$ws.Range("E42").Value = '  -2.46%  '
Set-TextValue 'D43' '8.012'
$ws.Range("E43").Value = '  -3.91%  '
Set-TextValue 'D44' '0.9995'
$ws.Range("E44").Value = '  -0.46%  '
Set-TextValue 'D45' '13.86'
$ws.Range("E45").Value = '  -2.73%  '
Set-TextValue 'D46' '0.6091'
$ws.Range("E46").Value = '  -1.74%  '
Set-TextValue 'D47' '3.830'
$ws.Range("E47").Value = '  +0.10%  '
Set-TextValue 'D48' '126.75'
$ws.Range("E48").Value = '  -3.18%  '
Set-TextValue 'D49' '2.011'
$ws.Range("E49").Value = '  -2.72%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue 'D50' '1.177'
$ws.Range("E50").Value = '  +3.35%  '
$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
Set-TextValue 'D51' '1.120'
$ws.Range("E51").Value = '  +6.41%  '
